$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): shift "prediction" and "rejection-f" left,
# and move "max" to the end.
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-7:
#  - Column C becomes a text duplicate of the species name (previously numeric)
#  - Column D stays the species name text (unchanged)
#  - Column E becomes a new numeric "rejection-f" score (previously species name text)

$ws.Range("C2").Value = "s__Clostridium_A leptum"
$ws.Range("E2").Value = 0.6626603308709718

$ws.Range("C3").Value = "s__Clostridium_A leptum"
$ws.Range("E3").Value = 0.6631453307174613

$ws.Range("C4").Value = "s__Clostridium_A leptum"
$ws.Range("E4").Value = 0.6737272617016222

$ws.Range("C5").Value = "s__Clostridium_A leptum"
$ws.Range("E5").Value = 0.6661637355733614

$ws.Range("C6").Value = "s__Clostridium_A leptum"
$ws.Range("E6").Value = 0.6453332538413822

$ws.Range("C7").Value = "s__Clostridium_A leptum"
$ws.Range("E7").Value = 0.671116898978992
